# Refresh the cryptocurrency price/volume snapshot (columns D = Price, E = Volume(1h))
# to match the latest scrape, as produced by the scheduled GitHub Actions job.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Values are written with a leading apostrophe so Excel keeps them as literal text
# (matching the workbook's original inline-string cells) instead of re-interpreting
# strings such as '1.00' or '0.790' as numbers and dropping the trailing zero.
$ws.Range("D2").Value = '''43.063.73'
$ws.Range("E2").Value = '''  -0.06%  '
$ws.Range("D3").Value = '''2.313.68'
$ws.Range("E3").Value = '''  +0.14%  '
$ws.Range("E4").Value = '''  +0.02%  '
$ws.Range("D5").Value = '''302.35'
$ws.Range("E5").Value = '''  -0.21%  '
$ws.Range("D6").Value = '''98.94'
$ws.Range("E6").Value = '''  -2.87%  '
$ws.Range("D7").Value = '''0.519'
$ws.Range("E7").Value = '''  +2.90%  '
$ws.Range("E9").Value = '''  -0.30%  '
$ws.Range("D10").Value = '''35.73'
$ws.Range("E10").Value = '''  -0.67%  '
$ws.Range("E11").Value = '''  -0.79%  '
$ws.Range("D12").Value = '''0.117'
$ws.Range("E12").Value = '''  -0.74%  '
$ws.Range("D13").Value = '''18.03'
$ws.Range("E13").Value = '''  +0.69%  '
$ws.Range("D14").Value = '''6.92'
$ws.Range("E14").Value = '''  -0.26%  '
$ws.Range("D15").Value = '''2.673.32'
$ws.Range("E15").Value = '''  -0.52%  '
$ws.Range("D16").Value = '''2.363.12'
$ws.Range("E16").Value = '''  +2.69%  '
$ws.Range("D17").Value = '''0.790'
$ws.Range("E17").Value = '''  -3.09%  '
$ws.Range("D18").Value = '''42.978.27'
$ws.Range("E18").Value = '''  -0.15%  '
$ws.Range("D19").Value = '''13.56'
$ws.Range("E19").Value = '''  +7.31%  '
$ws.Range("D20").Value = '''6.20'
$ws.Range("E20").Value = '''  -0.16%  '
$ws.Range("E21").Value = '''  +0.64%  '
$ws.Range("E22").Value = '''  +0.55%  '
$ws.Range("E23").Value = '''  +1.10%  '
$ws.Range("E24").Value = '''  -2.15%  '
$ws.Range("D25").Value = '''2.45'
$ws.Range("E25").Value = '''  -1.00%  '
$ws.Range("E26").Value = '''  -0.06%  '
$ws.Range("D27").Value = '''24.92'
$ws.Range("E27").Value = '''  +0.33%  '
$ws.Range("D28").Value = '''168.51'
$ws.Range("E28").Value = '''  +0.37%  '
$ws.Range("E29").Value = '''  -0.70%  '
$ws.Range("D30").Value = '''2.05'
$ws.Range("E30").Value = '''  -12.24%  '
$ws.Range("D31").Value = '''33.48'
$ws.Range("E31").Value = '''  -3.70%  '
$ws.Range("D32").Value = '''5.23'
$ws.Range("E32").Value = '''  +3.65%  '
$ws.Range("D33").Value = '''4.90'
$ws.Range("E33").Value = '''  +4.15%  '
$ws.Range("D34").Value = '''1.00'
$ws.Range("E34").Value = '''  -0.02%  '
$ws.Range("D35").Value = '''18.32'
$ws.Range("E35").Value = '''  +6.20%  '
$ws.Range("D37").Value = '''0.0695'
$ws.Range("E37").Value = '''  -0.36%  '
$ws.Range("E38").Value = '''  +0.18%  '
$ws.Range("E39").Value = '''  +0.45%  '
$ws.Range("E40").Value = '''  +0.94%  '
$ws.Range("E41").Value = '''  -2.24%  '
$ws.Range("D42").Value = '''2.000.87'
$ws.Range("E42").Value = '''  -0.05%  '
$ws.Range("E43").Value = '''  -0.23%  '
$ws.Range("D44").Value = '''2.17'
$ws.Range("E44").Value = '''  -5.75%  '
$ws.Range("E45").Value = '''  -1.80%  '
$ws.Range("D46").Value = '''17.49'
$ws.Range("E46").Value = '''  -1.56%  '
$ws.Range("E47").Value = '''  -2.07%  '
$ws.Range("D48").Value = '''54.89'
$ws.Range("E48").Value = '''  -2.12%  '
$ws.Range("D49").Value = '''2.539.19'
$ws.Range("E49").Value = '''  +0.73%  '
$ws.Range("D50").Value = '''74.24'
$ws.Range("E50").Value = '''  +5.42%  '
$ws.Range("E51").Value = '''  +0.43%  '
